# Green Line.xlsx — fill in the per-block Beacon Message codes (column J)
# Each block's beacon message changes from the generic placeholder
# "Beacon Msg Here" (or blank) to a unique code "G" + zero-padded block number,
# e.g. block 1 -> "G001", block 76 -> "G076", block 100 -> "G100".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$beaconCodes = @{
    2   = "G001"
    4   = "G003"
    9   = "G008"
    11  = "G010"
    13  = "G012"
    14  = "G013"
    16  = "G015"
    18  = "G017"
    22  = "G021"
    24  = "G023"
    29  = "G028"
    30  = "G029"
    31  = "G030"
    33  = "G032"
    39  = "G038"
    41  = "G040"
    48  = "G047"
    50  = "G049"
    57  = "G056"
    58  = "G057"
    59  = "G058"
    64  = "G063"
    65  = "G064"
    67  = "G066"
    73  = "G072"
    75  = "G074"
    77  = "G076"
    78  = "G077"
    79  = "G078"
    86  = "G085"
    87  = "G086"
    88  = "G087"
    90  = "G089"
    96  = "G095"
    98  = "G097"
    101 = "G100"
    102 = "G101"
    105 = "G104"
    107 = "G106"
    114 = "G113"
    116 = "G115"
    123 = "G122"
    125 = "G124"
    132 = "G131"
    134 = "G133"
    141 = "G140"
    143 = "G142"
    151 = "G150"
}

foreach ($rowNum in $beaconCodes.Keys) {
    $ws.Cells.Item($rowNum, 10).Value = $beaconCodes[$rowNum]
}

# Restore the view state roughly the way the author left it (scroll position
# and the active selection) — best-effort, cosmetic only.
$excel.ActiveWindow.ScrollRow = 114
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J131").Select()
